# Natmi following Dr Hou advice
# Update ligand/receptor-expressing cell counts (1 -> 3) and the
# dependent expression / specificity statistics for rows 2-17.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.379281
$ws.Range("H2").Value = 22.137843
$ws.Range("I2").Value = 0.2744121884499962
$ws.Range("J2").Value = 0.2744121884499961
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 87.667552
$ws.Range("N2").Value = 263.002656
$ws.Range("O2").Value = 0.3606416352150456
$ws.Range("P2").Value = 0.3606416352150456
$ws.Range("Q2").Value = 646.923500790112
$ws.Range("R2").Value = 5822.311507111008
$ws.Range("S2").Value = 0.09896446036554586
$ws.Range("T2").Value = 0.09896446036554585

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.379281
$ws.Range("H3").Value = 22.137843
$ws.Range("I3").Value = 0.2744121884499962
$ws.Range("J3").Value = 0.2744121884499961
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 72.97955566666666
$ws.Range("N3").Value = 218.938667
$ws.Range("O3").Value = 0.3002190170987564
$ws.Range("P3").Value = 0.3002190170987564
$ws.Range("Q3").Value = 538.5366485194755
$ws.Range("R3").Value = 4846.829836675281
$ws.Range("S3").Value = 0.08238375749637655
$ws.Range("T3").Value = 0.08238375749637654

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 7.379281
$ws.Range("H4").Value = 22.137843
$ws.Range("I4").Value = 0.2744121884499962
$ws.Range("J4").Value = 0.2744121884499961
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 41.37117366666666
$ws.Range("N4").Value = 124.113521
$ws.Range("O4").Value = 0.1701903085181653
$ws.Range("P4").Value = 0.1701903085181653
$ws.Range("Q4").Value = 305.2895157861336
$ws.Range("R4").Value = 2747.605642075203
$ws.Range("S4").Value = 0.04670229501344976
$ws.Range("T4").Value = 0.04670229501344975

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 7.379281
$ws.Range("H5").Value = 22.137843
$ws.Range("I5").Value = 0.2744121884499962
$ws.Range("J5").Value = 0.2744121884499961
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 41.06943633333334
$ws.Range("N5").Value = 123.208309
$ws.Range("O5").Value = 0.1689490391680327
$ws.Range("P5").Value = 0.1689490391680327
$ws.Range("Q5").Value = 303.0629112152764
$ws.Range("R5").Value = 2727.566200937487
$ws.Range("S5").Value = 0.04636167557462396
$ws.Range("T5").Value = 0.04636167557462396

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 13.29805733333333
$ws.Range("H6").Value = 39.894172
$ws.Range("I6").Value = 0.4945128143207339
$ws.Range("J6").Value = 0.4945128143207338
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 87.667552
$ws.Range("N6").Value = 263.002656
$ws.Range("O6").Value = 0.3606416352150456
$ws.Range("P6").Value = 0.3606416352150456
$ws.Range("Q6").Value = 1165.808132768981
$ws.Range("R6").Value = 10492.27319492083
$ws.Range("S6").Value = 0.1783419099914237
$ws.Range("T6").Value = 0.1783419099914237

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 13.29805733333333
$ws.Range("H7").Value = 39.894172
$ws.Range("I7").Value = 0.4945128143207339
$ws.Range("J7").Value = 0.4945128143207338
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 72.97955566666666
$ws.Range("N7").Value = 218.938667
$ws.Range("O7").Value = 0.3002190170987564
$ws.Range("P7").Value = 0.3002190170987564
$ws.Range("Q7").Value = 970.4863154165247
$ws.Range("R7").Value = 8734.376838748723
$ws.Range("S7").Value = 0.1484621510581106
$ws.Range("T7").Value = 0.1484621510581105

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 13.29805733333333
$ws.Range("H8").Value = 39.894172
$ws.Range("I8").Value = 0.4945128143207339
$ws.Range("J8").Value = 0.4945128143207338
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 41.37117366666666
$ws.Range("N8").Value = 124.113521
$ws.Range("O8").Value = 0.1701903085181653
$ws.Range("P8").Value = 0.1701903085181653
$ws.Range("Q8").Value = 550.1562393666235
$ws.Range("R8").Value = 4951.406154299611
$ws.Range("S8").Value = 0.0841612884354319
$ws.Range("T8").Value = 0.08416128843543187

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 13.29805733333333
$ws.Range("H9").Value = 39.894172
$ws.Range("I9").Value = 0.4945128143207339
$ws.Range("J9").Value = 0.4945128143207338
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 41.06943633333334
$ws.Range("N9").Value = 123.208309
$ws.Range("O9").Value = 0.1689490391680327
$ws.Range("P9").Value = 0.1689490391680327
$ws.Range("Q9").Value = 546.1437190083498
$ws.Range("R9").Value = 4915.293471075148
$ws.Range("S9").Value = 0.08354746483576773
$ws.Range("T9").Value = 0.08354746483576773

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.673314666666667
$ws.Range("H10").Value = 5.019944000000001
$ws.Range("I10").Value = 0.06222529534320158
$ws.Range("J10").Value = 0.06222529534320156
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 87.667552
$ws.Range("N10").Value = 263.002656
$ws.Range("O10").Value = 0.3606416352150456
$ws.Range("P10").Value = 0.3606416352150456
$ws.Range("Q10").Value = 146.6954005523627
$ws.Range("R10").Value = 1320.258604971264
$ws.Range("S10").Value = 0.02244103226431138
$ws.Range("T10").Value = 0.02244103226431137

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.673314666666667
$ws.Range("H11").Value = 5.019944000000001
$ws.Range("I11").Value = 0.06222529534320158
$ws.Range("J11").Value = 0.06222529534320156
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 72.97955566666666
$ws.Range("N11").Value = 218.938667
$ws.Range("O11").Value = 0.3002190170987564
$ws.Range("P11").Value = 0.3002190170987564
$ws.Range("Q11").Value = 122.1177608638498
$ws.Range("R11").Value = 1099.059847774648
$ws.Range("S11").Value = 0.0186812170066158
$ws.Range("T11").Value = 0.01868121700661579

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.673314666666667
$ws.Range("H12").Value = 5.019944000000001
$ws.Range("I12").Value = 0.06222529534320158
$ws.Range("J12").Value = 0.06222529534320156
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 41.37117366666666
$ws.Range("N12").Value = 124.113521
$ws.Range("O12").Value = 0.1701903085181653
$ws.Range("P12").Value = 0.1701903085181653
$ws.Range("Q12").Value = 69.22699167364712
$ws.Range("R12").Value = 623.042925062824
$ws.Range("S12").Value = 0.01059014221209343
$ws.Range("T12").Value = 0.01059014221209343

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.673314666666667
$ws.Range("H13").Value = 5.019944000000001
$ws.Range("I13").Value = 0.06222529534320158
$ws.Range("J13").Value = 0.06222529534320156
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 41.06943633333334
$ws.Range("N13").Value = 123.208309
$ws.Range("O13").Value = 0.1689490391680327
$ws.Range("P13").Value = 0.1689490391680327
$ws.Range("Q13").Value = 68.72209016829957
$ws.Range("R13").Value = 618.4988115146962
$ws.Range("S13").Value = 0.01051290386018096
$ws.Range("T13").Value = 0.01051290386018096

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.540576
$ws.Range("H14").Value = 13.621728
$ws.Range("I14").Value = 0.1688497018860685
$ws.Range("J14").Value = 0.1688497018860685
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 87.667552
$ws.Range("N14").Value = 263.002656
$ws.Range("O14").Value = 0.3606416352150456
$ws.Range("P14").Value = 0.3606416352150456
$ws.Range("Q14").Value = 398.061182589952
$ws.Range("R14").Value = 3582.550643309568
$ws.Range("S14").Value = 0.06089423259376472
$ws.Range("T14").Value = 0.06089423259376471

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.540576
$ws.Range("H15").Value = 13.621728
$ws.Range("I15").Value = 0.1688497018860685
$ws.Range("J15").Value = 0.1688497018860685
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 72.97955566666666
$ws.Range("N15").Value = 218.938667
$ws.Range("O15").Value = 0.3002190170987564
$ws.Range("P15").Value = 0.3002190170987564
$ws.Range("Q15").Value = 331.3692189507306
$ws.Range("R15").Value = 2982.322970556575
$ws.Range("S15").Value = 0.05069189153765352
$ws.Range("T15").Value = 0.05069189153765351

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.540576
$ws.Range("H16").Value = 13.621728
$ws.Range("I16").Value = 0.1688497018860685
$ws.Range("J16").Value = 0.1688497018860685
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 41.37117366666666
$ws.Range("N16").Value = 124.113521
$ws.Range("O16").Value = 0.1701903085181653
$ws.Range("P16").Value = 0.1701903085181653
$ws.Range("Q16").Value = 187.8489582426986
$ws.Range("R16").Value = 1690.640624184288
$ws.Range("S16").Value = 0.02873658285719024
$ws.Range("T16").Value = 0.02873658285719023

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 4.540576
$ws.Range("H17").Value = 13.621728
$ws.Range("I17").Value = 0.1688497018860685
$ws.Range("J17").Value = 0.1688497018860685
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 41.06943633333334
$ws.Range("N17").Value = 123.208309
$ws.Range("O17").Value = 0.1689490391680327
$ws.Range("P17").Value = 0.1689490391680327
$ws.Range("Q17").Value = 186.4788969486613
$ws.Range("R17").Value = 1678.310072537952
$ws.Range("S17").Value = 0.02852699489746003
$ws.Range("T17").Value = 0.02852699489746003
